$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 22
$ws.Range("C2").Value = 18

$ws.Range("A3").Value = "l2"
$ws.Range("C3").Value = 50

$ws.Range("A4").Value = "t1"
$ws.Range("B4").Value = "Televisor"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = $false

$ws.Rows("5:6").Delete()
